$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8 - pushes the existing rows 8-19 down to 9-20
$ws.Rows.Item(8).EntireRow.Insert()

# The A-L columns are identical across all the data rows (same market/product
# metadata), so copy them down from the row right below (old row 8, now row 9)
# into the freshly inserted row 8, preserving styles/number formats too.
$ws.Range("A9:L9").Copy()
$ws.Range("A8:L8").PasteSpecial()

# Fill in the new record's own data (row 8)
$ws.Range("D8").Value = "2023-01-27"
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 3000
$ws.Range("O8").Value = 3000
$ws.Range("P8").Value = 3000
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Provincia de Diguillín"
$ws.Range("S8").Value = 1500
$ws.Range("T8").Value = 2
